$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Change the language-separator syntax used by the import/export format ---
# "term|Lang;term|Lang" -> "term~~Lang|term~~Lang"
$ws.Range("G2").Value = "abrasión~~Espanhol|abrasion~~Inglês|εκτριβή~~Grego"
$ws.Range("G3").Value = "amarillamiento~~Espanhol|yellowing~~Inglês|κιτρίνισμα~~Grego"

# --- Taller header/data rows ---
$ws.Rows("1:4").RowHeight = 19.5

# --- Normalize font color from theme-based black to explicit RGB black ---
# Apply per contiguous block so each block keeps its own alignment / number
# format / border combination instead of merging into a single new style.
$ws.Range("A1:I1").Font.Color = 0
$ws.Range("A2:E2").Font.Color = 0
$ws.Range("G2:I2").Font.Color = 0
$ws.Range("A3:I3").Font.Color = 0
$ws.Range("A4:B4").Font.Color = 0
$ws.Range("J1").Font.Color = 0
$ws.Range("J2:J4").Font.Color = 0
